$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay plain text even when it looks
# like a number (e.g. "1.00", "0.0000247") by using the classic leading
# apostrophe, then clearing the resulting quote-prefix style so the cell
# keeps the workbook default (unstyled) formatting.
function Set-TextValue($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

$ws.Range("D2").Value = '69.080.68'
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").Value = '3.754.60'
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("E4").Value = '  +0.07%  '
Set-TextValue "D5" '602.18'
$ws.Range("E5").Value = '  +0.13%  '
Set-TextValue "D6" '166.36'
$ws.Range("E6").Value = '  -1.38%  '
$ws.Range("D7").Value = '3.752.63'
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +0.50%  '
$ws.Range("E10").Value = '  +4.14%  '
Set-TextValue "D11" '6.37'
$ws.Range("E11").Value = '  +1.12%  '
$ws.Range("E12").Value = '  -0.69%  '
Set-TextValue "D13" '37.68'
$ws.Range("E13").Value = '  -1.37%  '
Set-TextValue "D14" '0.0000247'
$ws.Range("E14").Value = '  +0.60%  '
$ws.Range("D15").Value = '4.386.87'
$ws.Range("E15").Value = '  +0.20%  '
$ws.Range("D16").Value = '3.755.40'
$ws.Range("E16").Value = '  -0.20%  '
$ws.Range("D17").Value = '69.107.62'
$ws.Range("E17").Value = '  +0.61%  '
$ws.Range("E18").Value = '  +1.54%  '
$ws.Range("E19").Value = '  +3.60%  '
$ws.Range("E20").Value = '  -1.02%  '
Set-TextValue "D21" '11.24'
$ws.Range("E21").Value = '  +3.08%  '
Set-TextValue "D22" '490.63'
$ws.Range("E22").Value = '  -0.93%  '
Set-TextValue "D23" '0.726'
$ws.Range("E23").Value = '  -0.48%  '
Set-TextValue "D24" '0.0000148'
$ws.Range("E24").Value = '  -1.46%  '
Set-TextValue "D25" '84.66'
$ws.Range("E25").Value = '  -0.95%  '
$ws.Range("E26").Value = '  -2.66%  '
Set-TextValue "D27" '12.26'
$ws.Range("E27").Value = '  -0.66%  '
Set-TextValue "D28" '10.05'
$ws.Range("E28").Value = '  -1.45%  '
$ws.Range("E29").Value = '  +0.12%  '
Set-TextValue "D30" '2.96'
$ws.Range("E30").Value = '  -0.61%  '
Set-TextValue "D31" '8.08'
$ws.Range("E31").Value = '  +1.86%  '
$ws.Range("E32").Value = '  -4.09%  '
Set-TextValue "D33" '31.67'
$ws.Range("E33").Value = '  -1.00%  '
$ws.Range("D34").Value = '3.904.13'
$ws.Range("E34").Value = '  +0.14%  '
$ws.Range("D35").Value = '3.708.43'
$ws.Range("E35").Value = '  +0.59%  '
$ws.Range("E36").Value = '  -0.53%  '
$ws.Range("E37").Value = '  +5.90%  '
Set-TextValue "D38" '5.92'
$ws.Range("E38").Value = '  +1.27%  '
$ws.Range("E39").Value = '  -1.12%  '
Set-TextValue "D40" '1.00'
$ws.Range("E40").Value = '  +0.03%  '
Set-TextValue "D41" '3.10'
$ws.Range("E41").Value = '  +7.84%  '
Set-TextValue "D42" '0.323'
$ws.Range("E42").Value = '  -0.02%  '
Set-TextValue "D43" '428.24'
$ws.Range("E43").Value = '  -3.28%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D44" '48.58'
$ws.Range("E44").Value = '  -0.67%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D45" '1.99'
$ws.Range("E45").Value = '  +0.81%  '
Set-TextValue "D46" '8.43'
$ws.Range("E46").Value = '  -0.39%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D48" '142.82'
$ws.Range("E48").Value = '  +0.94%  '
$ws.Range("B49").Value = 'Arweave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue "D49" '40.18'
$ws.Range("E49").Value = '  -0.54%  '
$ws.Range("D50").Value = '2.807.92'
$ws.Range("E50").Value = '  -0.56%  '
$ws.Range("E51").Value = '  +8.13%  '
